$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "name" column values so each row has a unique name
$ws.Range("B2").Value = "John Doe1"
$ws.Range("B3").Value = "John Doe2"
$ws.Range("B4").Value = "John Doe3"

# Add a new row 5, duplicating the pattern of rows 2-4, with a new employee
$ws.Range("A5").Value = "2024-12-20"
$ws.Range("B5").Value = "John Doe4"
$ws.Range("C5").Value = "123 Main St"
$ws.Range("D5").Value = 9876543210
$ws.Range("E5").Value = "john@example.com"
$ws.Range("F5").Value = "5 KW"
$ws.Range("G5").Value = "Company A"
$ws.Range("H5").Value = "Type 1"
$ws.Range("I5").Value = "Category 1"
$ws.Range("J5").Value = "Inverter Co A"
$ws.Range("K5").Value = "Category A"
$ws.Range("L5").Value = "Single Phase"
$ws.Range("M5").Value = "Type A"
$ws.Range("N5").Value = "Flat Roof"
$ws.Range("O5").Value = "Material A"
$ws.Range("P5").Value = "Material X"
$ws.Range("Q5").Value = "Earthing Type 1"
$ws.Range("R5").Value = "Wiring Type A"
$ws.Range("S5").Value = "DCDB Type 1"
$ws.Range("T5").Value = "ACDB Type 1"
$ws.Range("U5").Value = "Material 1"
$ws.Range("V5").Value = "System A"
$ws.Range("W5").Value = "Employee 1"
$ws.Range("X5").Value = "admin@example.com"

$ws.Hyperlinks.Add($ws.Range("X5"), "mailto:admin@example.com") | Out-Null

# Remove the now-unused trailing blank rows (previously rows 5 and 6 blank, row 7 blank)
$ws.Rows("6:7").Delete() | Out-Null
